$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new workout log rows (673-698) for the 25/03/2018 update.
$rows = @(
    @(672, 77, 43178, "March", 2018, "Monday", "Chest Fly", 105, 4, 8, "Chest"),
    @(673, 77, 43178, "March", 2018, "Monday", "Upright Rows", 30, 3, 12, "Shoulders"),
    @(674, 77, 43178, "March", 2018, "Monday", "Bicep Curl", 35, 4, 8, "Arms"),
    @(675, 77, 43178, "March", 2018, "Monday", "Seated Row", 80, 4, 8, "Back"),
    @(676, 78, 43180, "March", 2018, "Wednesday", "Incline Bench", 72.5, 4, 8, "Chest"),
    @(677, 78, 43180, "March", 2018, "Wednesday", "Bench Press", 80, 5, 8, "Chest"),
    @(678, 78, 43180, "March", 2018, "Wednesday", "Overhead Press", 40, 4, 8, "Shoulders"),
    @(679, 78, 43180, "March", 2018, "Wednesday", "Barbell Row", 80, 5, 5, "Back"),
    @(680, 79, 43182, "March", 2018, "Friday", "Shoulder Press", 25, 4, 8, "Shoulders"),
    @(681, 79, 43182, "March", 2018, "Friday", "Shoulder Shrug", 25, 4, 8, "Shoulders"),
    @(682, 79, 43182, "March", 2018, "Friday", "Tricep Pull down", 42.5, 4, 8, "Arms"),
    @(683, 79, 43182, "March", 2018, "Friday", "Hammer Curl", 17.5, 4, 8, "Arms"),
    @(684, 79, 43182, "March", 2018, "Friday", "Heel-taps", 0, 4, 15, "Core"),
    @(685, 79, 43182, "March", 2018, "Friday", "Raised leg circles", 0, 4, 12, "Core"),
    @(686, 79, 43182, "March", 2018, "Friday", "Scissors", 0, 4, 12, "Core"),
    @(687, 79, 43182, "March", 2018, "Friday", "Knee-Pull ins", 0, 4, 12, "Core"),
    @(688, 79, 43182, "March", 2018, "Friday", "Flitter Kicks", 0, 4, 20, "Core"),
    @(689, 80, 43184, "March", 2018, "Sunday", "Bench Press", 85, 5, 5, "Chest"),
    @(690, 80, 43184, "March", 2018, "Sunday", "Overhead Press", 50, 5, 5, "Shoulders"),
    @(691, 80, 43184, "March", 2018, "Sunday", "Barbell Row", 90, 5, 5, "Back"),
    @(692, 80, 43184, "March", 2018, "Sunday", "Pull-Ups", 93.8, 4, 5, "Back"),
    @(693, 80, 43184, "March", 2018, "Sunday", "Upright Rows", 30, 3, 12, "Shoulders"),
    @(694, 80, 43184, "March", 2018, "Sunday", "Bicep Curl", 30, 3, 12, "Arms"),
    @(695, 80, 43184, "March", 2018, "Sunday", "Plank", 0, 3, 30, "Core"),
    @(696, 80, 43184, "March", 2018, "Sunday", "Reverse Crunch", 0, 3, 12, "Core"),
    @(697, 80, 43184, "March", 2018, "Sunday", "Heel-taps", 0, 3, 15, "Core")
)

$startRow = 673
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Match the saved view state: scroll the frozen pane and move the active selection.
$excel.ActiveWindow.ScrollRow = 676
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C701").Select()
